# Stand up weeks 8, 9 y 10
# Fill in the "Paula Andrea Taborda Jaramillo" row block (rows 10-12) with
# her Monday (C) / Tuesday (D) stand-up answers for week 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "Terminé de coregir el caso de uso"
$ws.Range("D10").Value = "Nada"
$ws.Range("C11").Value = "Nada"
$ws.Range("C12").Value = "Otras materias"
$ws.Range("D11").Value = "Asistir a la reunión con los compañeros"
$ws.Range("D12").Value = "Ninguna"

$ws.Range("D12").Select()
